$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.073.67"
$ws.Range("E2").Value = "  +5.99%  "
$ws.Range("D3").Value = "2.733.32"
$ws.Range("E3").Value = "  +4.56%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'585.85"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").Value = "'151.19"
$ws.Range("E6").Value = "  +5.31%  "
$ws.Range("D7").Value = "'0.995"
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("D8").Value = "'0.611"
$ws.Range("E8").Value = "  +2.52%  "
$ws.Range("D9").Value = "2.767.13"
$ws.Range("E9").Value = "  +5.56%  "
$ws.Range("D10").Value = "'6.77"
$ws.Range("E10").Value = "  +3.97%  "
$ws.Range("E11").Value = "  +8.16%  "
$ws.Range("E12").Value = "  +4.31%  "
$ws.Range("E13").Value = "  +1.55%  "
$ws.Range("D14").Value = "3.221.18"
$ws.Range("E14").Value = "  +4.78%  "
$ws.Range("D15").Value = "'26.72"
$ws.Range("E15").Value = "  +8.11%  "
$ws.Range("D16").Value = "63.929.20"
$ws.Range("E16").Value = "  +5.79%  "
$ws.Range("E17").Value = "  +8.04%  "
$ws.Range("D18").Value = "2.757.25"
$ws.Range("E18").Value = "  +5.34%  "
$ws.Range("D19").Value = "'12.02"
$ws.Range("E19").Value = "  +5.76%  "
$ws.Range("E20").Value = "  +5.30%  "
$ws.Range("D21").Value = "'367.35"
$ws.Range("E21").Value = "  +5.89%  "
$ws.Range("D22").Value = "'7.02"
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "'0.534"
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("D25").Value = "'66.02"
$ws.Range("E25").Value = "  +3.76%  "
$ws.Range("D26").Value = "'0.168"
$ws.Range("E26").Value = "  +4.42%  "
$ws.Range("E27").Value = "  +8.50%  "
$ws.Range("D28").Value = "'0.994"
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("D29").Value = "0.0₃0880"
$ws.Range("E29").Value = "  +10.02%  "
$ws.Range("E30").Value = "  +6.83%  "
$ws.Range("E31").Value = "  +10.70%  "
$ws.Range("D32").Value = "'170.39"
$ws.Range("E32").Value = "  +0.95%  "
$ws.Range("D33").Value = "'1.20"
$ws.Range("E33").Value = "  +18.81%  "
$ws.Range("D34").Value = "'0.996"
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("D35").Value = "'20.62"
$ws.Range("E35").Value = "  +5.88%  "
$ws.Range("D36").Value = "'4.78"
$ws.Range("E36").Value = "  +11.55%  "
$ws.Range("D37").Value = "'1.44"
$ws.Range("E37").Value = "  +10.08%  "
$ws.Range("D38").Value = "'1.82"
$ws.Range("E38").Value = "  +10.13%  "
$ws.Range("E39").Value = "  +19.83%  "
$ws.Range("D40").Value = "'354.09"
$ws.Range("E40").Value = "  +10.67%  "
$ws.Range("D41").Value = "'4.26"
$ws.Range("E41").Value = "  +9.15%  "
$ws.Range("D42").Value = "'39.38"
$ws.Range("E42").Value = "  +2.61%  "
$ws.Range("D43").Value = "'5.68"
$ws.Range("E43").Value = "  +12.84%  "
$ws.Range("D44").Value = "'22.36"
$ws.Range("E44").Value = "  +12.09%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'143.47"
$ws.Range("E45").Value = "  +5.45%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'22.06"
$ws.Range("E46").Value = "  +9.75%  "
$ws.Range("D47").Value = "'0.0592"
$ws.Range("E47").Value = "  +7.59%  "
$ws.Range("D48").Value = "'0.645"
$ws.Range("E48").Value = "  +5.60%  "
$ws.Range("E49").Value = "  +7.44%  "
$ws.Range("D50").Value = "'0.102"
$ws.Range("E50").Value = "  +2.25%  "
$ws.Range("D51").Value = "2.171.09"
$ws.Range("E51").Value = "  +7.00%  "
